# STAR_template.docx minor text updates ("minor updates to testing reports")
#
# 1. "Your child has just taken a STAR Early Literacy assessment on our
#     school" -> "Your child has taken a STAR Early Literacy
#     assessment(s) on our school"          (drop "just ", add "(s)")
# 2. The Word "_GoBack" bookmark (marks the location of the author's last
#    edit) moves from its old spot (inside "...received on the STAR Early
#    Literacy assessment...") to a new spot inside "Your child has taken
#    this assessment" (after "Your ch").  Re-adding a bookmark with the
#    same name moves it, so we just add it at the new location and then
#    tidy up the text run at its old location.

$d = $word.ActiveDocument

# --- Hunk 1: "just taken" -> "taken" ---------------------------------
$r = $d.Content
$r.Find.Execute("just taken", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "taken", 2) | Out-Null

# --- Hunk 1: "assessment on our school" -> "assessment(s) on our school"
$r = $d.Content
$r.Find.Execute("assessment on our school", $true, $false, $false, $false, `
                 $false, $true, 1, $false, "assessment(s) on our school", 2) | Out-Null

# --- Hunk 2: move the _GoBack bookmark into "Your child has taken this
#     assessment " -- right after "Your ch" -------------------------
$r = $d.Content
$r.Find.Execute("Your child has taken this assessment ", $true, $false, `
                 $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmStart = $r.Start + 7
$bmRange = $d.Range($bmStart, $bmStart)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# --- Hunk 3: tidy up the text run left behind at the bookmark's old
#     location ("...received on th" | "e STAR Early Literacy...") so it
#     reads as a single run again, now that the bookmark has moved away.
$r = $d.Content
$r.Find.Execute("received on th", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "received on th", 2) | Out-Null
